$wb = $excel.ActiveWorkbook

# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet -- this pushes the existing "Late" / "Heading" / "Outstanding"
# columns one slot to the right (N->O, O->P, P->Q) and leaves the new N
# column empty.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
[void]$wsSchedule.Columns("N").Insert()

# Make the "Repayment Schedule" sheet the active tab / selected sheet, and
# restore the cursor to cell S8 on it (matches the author's last selection
# when they saved the workbook).
[void]$wsSchedule.Activate()
[void]$wsSchedule.Range("S8").Select()
